# CU-10 EliminarProveedor - CU-10_Descripcion.docx edits
# Applies the corrections described in the commit:
#   - wraps "ConfirmationView"/"ErrorView" occurrences with spellcheck
#     proofErr markers (splitting the host run in three)
#   - drops the stray "  (FA-01)" annotation on the normal flow
#   - adds a new "El FA-01 puede ocurrir..." note at the end of the
#     normal flow
#   - rewords the EX-01 exception trigger
#   - tidies a couple of runs that were split/merged unnecessarily
#
# Strategy: locate each target paragraph with Find.Execute against a
# short, unique anchor string, then replace that whole paragraph's
# content (via Range.InsertXML on the paragraph's Range, which keeps
# the paragraph's own pPr/pStyle/numPr intact when supplied again) with
# the corrected run/proofErr layout.

$d = $word.ActiveDocument
$pkgHeader = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>'
$pkgFooter = '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

function Replace-Paragraph($anchorText, $newParagraphsXml) {
    $rng = $d.Content
    $found = $rng.Find.Execute($anchorText)
    if (-not $found) {
        throw "Anchor not found: $anchorText"
    }
    $prng = $rng.Paragraphs(1).Range
    $xml = $pkgHeader + $newParagraphsXml + $pkgFooter
    $prng.InsertXML($xml)
}

# --- Hunk 1 -----------------------------------------------------------
# "la ventana ConfirmationView" -> "la ventana " + proofErr(ConfirmationView)
Replace-Paragraph "la ventana ConfirmationView" (
    '<w:p><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:jc w:val="both"/></w:pPr>' +
    '<w:r><w:t xml:space="preserve">El sistema muestra </w:t></w:r>' +
    '<w:r><w:t xml:space="preserve">la ventana </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>ConfirmationView</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t xml:space="preserve"> con el mensaje &#8220;</w:t></w:r>' +
    '<w:r><w:t>&#191;Est&#225; seguro de que desea eliminar a este proveedor? Esta acci&#243;n no se puede deshacer</w:t></w:r>' +
    '<w:r><w:t>&#8221; y un bot&#243;n &#8220;Aceptar&#8221;</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> y un bot&#243;n &#8220;Cancelar&#8221;</w:t></w:r>' +
    '</w:p>'
)

# --- Hunk 2 -----------------------------------------------------------
# Drop the "  (FA-01)" run entirely.
Replace-Paragraph "(FA-01)" (
    '<w:p><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:jc w:val="both"/></w:pPr>' +
    '<w:r><w:t>El actor hace clic en el bot&#243;n &#8220;Aceptar&#8221;.</w:t></w:r>' +
    '</w:p>'
)

# --- Hunk 3 -------------------------------------------------------------
# "ConfirmationView. " -> proofErr(ConfirmationView) + ". "
Replace-Paragraph "cambia el estado del PROVEEDOR" (
    '<w:p><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:jc w:val="both"/></w:pPr>' +
    '<w:r><w:t>El sistema</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
    '<w:r><w:t>cambia el estado del PROVEEDOR dentro de la base de datos</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> y </w:t></w:r>' +
    '<w:r><w:t xml:space="preserve">cierra </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>ConfirmationView</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t xml:space="preserve">. </w:t></w:r>' +
    '<w:r><w:t>(EX-01)</w:t></w:r>' +
    '</w:p>'
)

# --- Hunk 4 -------------------------------------------------------------
# Keep "Fin del caso de uso." but append an empty paragraph and a new
# "El FA-01 puede ocurrir..." paragraph right after it (still inside the
# same table cell, i.e. same Find/Paragraphs(1) target).
Replace-Paragraph "Fin del caso de uso." (
    '<w:p><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:jc w:val="both"/></w:pPr>' +
    '<w:r><w:t>Fin del caso de uso.</w:t></w:r>' +
    '</w:p>' +
    '<w:p><w:pPr><w:jc w:val="both"/></w:pPr></w:p>' +
    '<w:p><w:pPr><w:jc w:val="both"/></w:pPr>' +
    '<w:r><w:t>El FA-01 puede ocurrir en cualquier momento del CU</w:t></w:r>' +
    '</w:p>'
)

# --- Hunk 5 -------------------------------------------------------------
# "cierra ConfirmationView." (Flujo alterno) -> "cierra " + proofErr(ConfirmationView) + "."
Replace-Paragraph "El sistema cierra ConfirmationView" (
    '<w:p><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr><w:jc w:val="both"/></w:pPr>' +
    '<w:r><w:t xml:space="preserve">El sistema </w:t></w:r>' +
    '<w:r><w:t xml:space="preserve">cierra </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>ConfirmationView</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t>.</w:t></w:r>' +
    '</w:p>'
)

# --- Hunk 6 -------------------------------------------------------------
# "EX-01 No hay conexión con la base de datos" -> "EX-01 No hay conexión " + "a la red"
Replace-Paragraph "EX-01 No hay conexi" (
    '<w:p><w:pPr><w:jc w:val="both"/></w:pPr>' +
    '<w:r><w:t xml:space="preserve">EX-01 No hay conexi&#243;n </w:t></w:r>' +
    '<w:r><w:t>a la red</w:t></w:r>' +
    '</w:p>'
)

# --- Hunk 7 -------------------------------------------------------------
# "la ventana ErrorView con el mensaje " -> "la ventana " + proofErr(ErrorView) + " con el mensaje "
Replace-Paragraph "El sistema muestra la ventana ErrorView" (
    '<w:p><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="4"/></w:numPr><w:jc w:val="both"/></w:pPr>' +
    '<w:r><w:t xml:space="preserve">El sistema muestra la ventana </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>ErrorView</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t xml:space="preserve"> con el mensaje </w:t></w:r>' +
    '<w:r><w:t>&#8220;</w:t></w:r>' +
    '<w:r><w:t>No se pudo conectar a la red del supermercado, int&#233;ntelo de nuevo m&#225;s tarde</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve">&#8221; </w:t></w:r>' +
    '<w:r><w:t>y un bot&#243;n &#8220;Aceptar&#8221;.</w:t></w:r>' +
    '</w:p>'
)

# --- Hunk 8 -------------------------------------------------------------
# "El sistema cierra la ventana ErrorView" -> "...la ventana " + proofErr(ErrorView)
Replace-Paragraph "El sistema cierra la ventana ErrorView" (
    '<w:p><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="4"/></w:numPr><w:jc w:val="both"/></w:pPr>' +
    '<w:r><w:t xml:space="preserve">El sistema cierra la ventana </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>ErrorView</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t>.</w:t></w:r>' +
    '</w:p>'
)

# --- Hunk 9 -------------------------------------------------------------
# "Fin del caso de uso" + "." -> single run "Fin del caso de uso."
Replace-Paragraph "Fin del caso de uso</w:t" (
    '<w:p><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="4"/></w:numPr></w:pPr>' +
    '<w:r><w:t>Fin del caso de uso.</w:t></w:r>' +
    '</w:p>'
)

# --- Hunk 10 --------------------------------------------------------------
# "El estado del PROVEEDOR cambia a Eliminado" -> "...cambia a " + "Eliminado"
Replace-Paragraph "El estado del PROVEEDOR cambia a" (
    '<w:p>' +
    '<w:r><w:t xml:space="preserve">POS-01 </w:t></w:r>' +
    '<w:r><w:t xml:space="preserve">El estado del PROVEEDOR cambia a </w:t></w:r>' +
    '<w:r><w:t>Eliminado</w:t></w:r>' +
    '</w:p>'
)

Write-Host "Done"
